$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.005.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.101.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.89%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "349.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.61%  "

$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5147"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.98%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.15"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08971"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.094.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.194"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.723"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001144"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06662"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.211"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.103.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.337"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.343.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.560"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.171"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1059"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.643"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("E34").Value = "  -2.40%  "

$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.135"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02568"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06776"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2263"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.339"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.46%  "

$ws.Range("E42").Value = "  -1.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6796"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6422"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.278"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000363"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.638"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("E49").Value = "  -3.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.05"
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07216"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.33%  "
